# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN.. right by one)
#  - change Week labels from zero-padded (W01..W09) to unpadded (W1..W9)
#  - populate the new Week_Start_Date column with the week's start date (as text)
#  - recompute Amazon P70/P80/P90 forecast values for the shifted columns
#  - store is_holiday_week as a boolean (FALSE) instead of a number

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the old column B (ASIN), shifting C.. -> D..
$ws.Columns.Item(2).Insert()

# New header cell for the inserted column
$ws.Range("B1").Value = "Week_Start_Date"

$weekLabels = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # A: Week label, unpadded (W01 -> W1, etc.)
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]

    # B: new Week_Start_Date column - force text so Excel doesn't coerce to a date serial
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $weekStartDates[$i]

    # F (Amazon P70 Forecast, shifted from old E): now 0 instead of 1
    $ws.Cells.Item($row, 6).Value = 0

    # G (Amazon P80 Forecast, shifted from old F): stays 1
    $ws.Cells.Item($row, 7).Value = 1

    # H (Amazon P90 Forecast, shifted from old G): now 1 (numeric, not the title text)
    $ws.Cells.Item($row, 8).Value = 1

    # J (is_holiday_week, shifted from old I): boolean FALSE instead of numeric 0
    $ws.Cells.Item($row, 10).Value = $false
}
